$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps displaying numeric-looking values as literal text
# (matches the source data which stores prices/changes as plain strings).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2: Bitcoin
$ws.Range("D2").Value = "28.754.99"
$ws.Range("E2").Value = "  +2.26%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.814.13"
$ws.Range("E3").Value = "  -0.29%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5: BNB
$ws.Range("D5").Value = "329.25"
$ws.Range("E5").Value = "  -2.53%  "

# Row 6: USDC
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.33%  "

# Row 7: XRP
$ws.Range("D7").Value = "0.4379"
$ws.Range("E7").Value = "  +1.93%  "

# Row 8: Cardano
$ws.Range("D8").Value = "0.3793"
$ws.Range("E8").Value = "  +7.83%  "

# Row 9: OKB
$ws.Range("D9").Value = "44.49"
$ws.Range("E9").Value = "  -2.35%  "

# Row 10: Dogecoin
$ws.Range("D10").Value = "0.07788"
$ws.Range("E10").Value = "  +4.51%  "

# Row 11: Polygon
$ws.Range("D11").Value = "1.149"
$ws.Range("E11").Value = "  -0.43%  "

# Row 12: Solana
$ws.Range("D12").Value = "22.75"
$ws.Range("E12").Value = "  -1.30%  "

# Row 13: BinanceUSD
$ws.Range("D13").Value = "0.9986"
$ws.Range("E13").Value = "  -0.13%  "

# Row 14: Chainlink
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "7.654"
$ws.Range("E14").Value = "  +5.01%  "

# Row 15: Polkadot
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "6.331"
$ws.Range("E15").Value = "  +1.06%  "

# Row 16: WrappedEther
$ws.Range("D16").Value = "1.810.11"
$ws.Range("E16").Value = "  -0.30%  "

# Row 17: ShibaInu
$ws.Range("D17").Value = "0.00001100"
$ws.Range("E17").Value = "  +1.25%  "

# Row 18: TRON
$ws.Range("D18").Value = "0.06762"
$ws.Range("E18").Value = "  +1.16%  "

# Row 19: Litecoin
$ws.Range("D19").Value = "81.57"
$ws.Range("E19").Value = "  -0.60%  "

# Row 20: Dai
$ws.Range("D20").Value = "0.9998"
$ws.Range("E20").Value = "  -0.01%  "

# Row 21: Avalanche
$ws.Range("D21").Value = "17.75"
$ws.Range("E21").Value = "  +2.74%  "

# Row 22: Uniswap
$ws.Range("D22").Value = "6.321"
$ws.Range("E22").Value = "  -2.84%  "

# Row 23: WrappedBTC
$ws.Range("D23").Value = "28.688.38"
$ws.Range("E23").Value = "  +1.96%  "

# Row 24: Cosmos
$ws.Range("D24").Value = "11.85"
$ws.Range("E24").Value = "  -1.42%  "

# Row 25: Toncoin
$ws.Range("D25").Value = "2.442"
$ws.Range("E25").Value = "  +1.93%  "

# Row 26: EthereumClassic
$ws.Range("D26").Value = "20.69"
$ws.Range("E26").Value = "  -0.24%  "

# Row 27: Monero
$ws.Range("D27").Value = "153.26"
$ws.Range("E27").Value = "  -1.96%  "

# Row 28: LidoDAOToken
$ws.Range("D28").Value = "2.373"
$ws.Range("E28").Value = "  -4.94%  "

# Row 29: WrappedliquidstakedEther2.0
$ws.Range("D29").Value = "2.017.92"
$ws.Range("E29").Value = "  -0.25%  "

# Row 30: ImmutableX
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "1.282"
$ws.Range("E30").Value = "  -1.80%  "

# Row 31: BitcoinCash
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "132.58"
$ws.Range("E31").Value = "  +0.03%  "

# Row 32: HuobiToken
$ws.Range("D32").Value = "3.977"
$ws.Range("E32").Value = "  -1.92%  "

# Row 33: Filecoin
$ws.Range("D33").Value = "5.857"
$ws.Range("E33").Value = "  -1.96%  "

# Row 34: Stellar
$ws.Range("D34").Value = "0.09304"
$ws.Range("E34").Value = "  +0.68%  "

# Row 35: Algorand
$ws.Range("D35").Value = "0.2270"
$ws.Range("E35").Value = "  +4.73%  "

# Row 36: Aptos
$ws.Range("D36").Value = "12.26"
$ws.Range("E36").Value = "  -0.98%  "

# Row 37: Hedera
$ws.Range("D37").Value = "0.06392"
$ws.Range("E37").Value = "  +2.14%  "

# Row 38: VeChain
$ws.Range("D38").Value = "0.02348"
$ws.Range("E38").Value = "  -0.99%  "

# Row 39: InternetComputer(DFINITY)
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "5.241"
$ws.Range("E39").Value = "  -0.21%  "

# Row 40: TheSandbox
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.6664"
$ws.Range("E40").Value = "  -1.36%  "

# Row 41: TrustWalletToken
$ws.Range("D41").Value = "1.206"
$ws.Range("E41").Value = "  -1.06%  "

# Row 42: FraxShare
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "8.169"
$ws.Range("E42").Value = "  -0.81%  "

# Row 43: WEMIXTOKEN
$ws.Range("B43").Value = "WEMIXTOKEN"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "1.450"
$ws.Range("E43").Value = "  -2.60%  "

# Row 44: Frax
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.23%  "

# Row 45: EnergySwap
$ws.Range("D45").Value = "14.00"
$ws.Range("E45").Value = "  -0.06%  "

# Row 46: Decentraland
$ws.Range("D46").Value = "0.6113"
$ws.Range("E46").Value = "  -0.48%  "

# Row 47: PancakeSwap
$ws.Range("D47").Value = "3.814"
$ws.Range("E47").Value = "  -1.63%  "

# Row 48: Quant
$ws.Range("D48").Value = "129.09"
$ws.Range("E48").Value = "  +0.17%  "

# Row 49: NEARProtocol
$ws.Range("D49").Value = "2.058"
$ws.Range("E49").Value = "  +0.43%  "

# Row 50: Cronos
$ws.Range("D50").Value = "0.07104"
$ws.Range("E50").Value = "  -0.04%  "

# Row 51: EOS
$ws.Range("D51").Value = "1.150"
$ws.Range("E51").Value = "  -2.51%  "

